$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A80").Value = "JZT1M4"
$ws.Range("B80").Value = "Motor para contador de billetes"
$ws.Range("C80").Value = "9V 110Rpm"
$ws.Range("D80").Value = 0
$ws.Range("E80").Value = 180000
$ws.Range("F80").Value = 2
$ws.Range("G80").Value = 2
$ws.Range("H80").Formula = "=(E80-D80)*G80"
$ws.Range("I80").Formula = "=D80*F80"
$ws.Range("J80").Value = 0
